$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-04-19 Friday" "2024-04-20 Saturday"
Replace-Text "139×6=" "226×8="
Replace-Text "692×3=" "323×9="
Replace-Text "524×4=" "847×2="
Replace-Text "217×3=" "441×7="
Replace-Text "779×4=" "804×5="
Replace-Text "464×2=" "814×6="
Replace-Text "955×9=" "875×9="
Replace-Text "458×4=" "310×6="
Replace-Text "642×4=" "844×8="
Replace-Text "193×2=" "161×7="
Replace-Text "884×2=" "220×8="
Replace-Text "238×2=" "778×9="
Replace-Text "229×2=" "817×6="
Replace-Text "422×6=" "290×9="
Replace-Text "495×7=" "787×9="
Replace-Text "411×8=" "628×2="
Replace-Text "830×4=" "367×9="
Replace-Text "713×4=" "916×6="
Replace-Text "964×7=" "485×7="
Replace-Text "962×4=" "759×8="
Replace-Text "775×3=" "530×4="
Replace-Text "948×8=" "122×3="
Replace-Text "459×7=" "399×7="
Replace-Text "955×6=" "128×2="
Replace-Text "852×6=" "276×2="

Write-Output "Done"
